$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPbES")
$ws.Activate()

# Row 2 (hard coal): guaranteed-dispatch flag changes from 1 to 0 across B2:AE2
$ws.Range("B2:AE2").Value = 0

# Row 19 (natural gas combined cycle): guaranteed-dispatch flag changes from 0 to 1 across B19:AE19
$ws.Range("B19:AE19").Value = 1

# Update the active selection to match the saved workbook state
$ws.Range("B2:AE2").Select()
